# Projeto finalizado com persistência de dados quando executado em ambiente local.
#
# Appends a new "motivos" log entry (row 7) - a rerun of the same
# Usina 1 / Inversor / 331 / Bom check already logged in row 4 - and keeps
# the "data" column's date format consistent between the last two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 7: same shape as row 4 (Usina 1 / Inversor / 331 / Bom) ---
$ws.Cells.Item(7, 1).Value = "Usina 1"
$ws.Cells.Item(7, 2).Value = "Inversor"
$ws.Cells.Item(7, 4).Value = "Bom"
# "detalhamento" (column E) is left blank for this entry, same as row 4.

# "data" column: same date serial as the row above it.
$ws.Cells.Item(7, 6).Value = 45677

# Sync F7's date format to F6's, then refresh F6 onto the same format so
# both end up sharing a single, consistent date style.
$ws.Cells.Item(7, 6).NumberFormat = "yyyy-mm-dd"
$ws.Cells.Item(6, 6).NumberFormat = "yyyy-mm-dd"

# "id_equipamento" is kept as text (matches every other cell in column C),
# so force a text format before assigning the numeric-looking value, then
# drop back to the default style so no stray number-format lingers on the
# cell itself.
$ws.Cells.Item(7, 3).NumberFormat = "@"
$ws.Cells.Item(7, 3).Value = "331"
$ws.Cells.Item(7, 3).Style = "Normal"

Write-Output "Row 7 appended (Usina 1 / Inversor / 331 / Bom); F6/F7 date format synced."
